$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new rows at position 29 (shifts old rows 29-79 down to 31-81)
$ws.Rows("29:30").Insert()

# New data for row 29 (Ají, Inferno, Primera - 12kg box, serial date 44428)
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44428
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112021
$ws.Range("G29").Value = "Ají"
$ws.Range("H29").Value = "Inferno"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 700
$ws.Range("K29").Value = 36000
$ws.Range("L29").Value = 37000
$ws.Range("M29").Value = 36500
$ws.Range("N29").Value = "$/caja 12 kilos"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 3042
$ws.Range("Q29").Value = 12
$ws.Range("R29").Value = "Hortaliza"

# New data for row 30 (Ají, Inferno, Segunda - 12kg box, serial date 44428)
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44428
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112021
$ws.Range("G30").Value = "Ají"
$ws.Range("H30").Value = "Inferno"
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 29000
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = 29500
$ws.Range("N30").Value = "$/caja 12 kilos"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 2458
$ws.Range("Q30").Value = 12
$ws.Range("R30").Value = "Hortaliza"

# Step 2: insert one new row at position 55 (shifts rows down by 1 further)
$ws.Rows("55:55").Insert()

# New data for row 55 (Ají, Inferno, Primera - 12kg box, serial date 44435)
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44435
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112021
$ws.Range("G55").Value = "Ají"
$ws.Range("H55").Value = "Inferno"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 1840
$ws.Range("K55").Value = 34000
$ws.Range("L55").Value = 36000
$ws.Range("M55").Value = 35120
$ws.Range("N55").Value = "$/caja 12 kilos"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 2927
$ws.Range("Q55").Value = 12
$ws.Range("R55").Value = "Hortaliza"

# Step 3: insert one new row at position 69 (shifts rows down by 1 further)
$ws.Rows("69:69").Insert()

# New data for row 69 (Ají, Inferno, Primera - 12kg box, serial date 44433)
$ws.Range("A69").Value = 8
$ws.Range("B69").Value = "Terminal La Palmera de La Serena"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 44433
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Inferno"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 600
$ws.Range("K69").Value = 35000
$ws.Range("L69").Value = 36000
$ws.Range("M69").Value = 35500
$ws.Range("N69").Value = "$/caja 12 kilos"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 2958
$ws.Range("Q69").Value = 12
$ws.Range("R69").Value = "Hortaliza"

# Step 4: append one new row at the very end (position 84)
$ws.Rows("84:84").Insert()

# New data for row 84 (Ají, Inferno, Primera - 12kg box, serial date 44432)
$ws.Range("A84").Value = 8
$ws.Range("B84").Value = "Terminal La Palmera de La Serena"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44432
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 100112021
$ws.Range("G84").Value = "Ají"
$ws.Range("H84").Value = "Inferno"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 540
$ws.Range("K84").Value = 35000
$ws.Range("L84").Value = 36000
$ws.Range("M84").Value = 35500
$ws.Range("N84").Value = "$/caja 12 kilos"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 2958
$ws.Range("Q84").Value = 12
$ws.Range("R84").Value = "Hortaliza"

